$wb = $excel.ActiveWorkbook

# --- three_line: add rows 199-204 ---
$ws = $wb.Worksheets.Item("three_line")
$ws.Range("A199").Value = 45439.55208333334
$ws.Range("A199").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B199").Value = '11-06-2024 10:15:00'
$ws.Range("C199").Value = 'hour'
$ws.Range("D199").Value = 'TATASTEEL.NS'
$ws.Range("E199").Value = 45408.38541666666
$ws.Range("E199").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F199").Value = 170.6999969482422
$ws.Range("G199").Value = 45434.38541666666
$ws.Range("G199").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H199").Value = 175.4499969482422
$ws.Range("I199").Value = 45436.38541666666
$ws.Range("I199").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J199").Value = 177.5500030517578
$ws.Range("K199").Value = 'High'
$ws.Range("L199").Value = '11/06/2024 05:47:02'

$ws.Range("A200").Value = 45441.55208333334
$ws.Range("A200").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B200").Value = '11-06-2024 10:15:00'
$ws.Range("C200").Value = 'hour'
$ws.Range("D200").Value = 'TATASTEEL.NS'
$ws.Range("E200").Value = 45415.38541666666
$ws.Range("E200").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F200").Value = 170.75
$ws.Range("G200").Value = 45434.38541666666
$ws.Range("G200").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H200").Value = 175.4499969482422
$ws.Range("I200").Value = 45440.38541666666
$ws.Range("I200").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J200").Value = 177.5
$ws.Range("K200").Value = 'High'
$ws.Range("L200").Value = '11/06/2024 05:47:02'

$ws.Range("A201").Value = 45441.55208333334
$ws.Range("A201").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B201").Value = '11-06-2024 10:15:00'
$ws.Range("C201").Value = 'hour'
$ws.Range("D201").Value = 'TATASTEEL.NS'
$ws.Range("E201").Value = 45415.38541666666
$ws.Range("E201").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F201").Value = 170.75
$ws.Range("G201").Value = 45436.38541666666
$ws.Range("G201").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H201").Value = 177.5500030517578
$ws.Range("I201").Value = 45440.38541666666
$ws.Range("I201").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J201").Value = 177.5
$ws.Range("K201").Value = 'High'
$ws.Range("L201").Value = '11/06/2024 05:47:02'

$ws.Range("A202").Value = 45447.55208333334
$ws.Range("A202").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B202").Value = '11-06-2024 10:15:00'
$ws.Range("C202").Value = 'hour'
$ws.Range("D202").Value = 'MHRIL.NS'
$ws.Range("E202").Value = 45433.38541666666
$ws.Range("E202").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F202").Value = 412.7999877929688
$ws.Range("G202").Value = 45435.63541666666
$ws.Range("G202").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H202").Value = 412
$ws.Range("I202").Value = 45442.46875
$ws.Range("I202").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J202").Value = 409.75
$ws.Range("K202").Value = 'High'
$ws.Range("L202").Value = '11/06/2024 05:47:02'

$ws.Range("A203").Value = 45447.55208333334
$ws.Range("A203").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B203").Value = '11-06-2024 10:15:00'
$ws.Range("C203").Value = 'hour'
$ws.Range("D203").Value = 'MHRIL.NS'
$ws.Range("E203").Value = 45433.38541666666
$ws.Range("E203").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F203").Value = 412.7999877929688
$ws.Range("G203").Value = 45442.42708333334
$ws.Range("G203").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H203").Value = 409.75
$ws.Range("I203").Value = 45442.46875
$ws.Range("I203").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J203").Value = 409.75
$ws.Range("K203").Value = 'High'
$ws.Range("L203").Value = '11/06/2024 05:47:02'

$ws.Range("A204").Value = 45436.55208333334
$ws.Range("A204").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B204").Value = '11-06-2024 10:15:00'
$ws.Range("C204").Value = 'hour'
$ws.Range("D204").Value = 'RAJSREESUG.NS'
$ws.Range("E204").Value = 45420.55208333334
$ws.Range("E204").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F204").Value = 65.25
$ws.Range("G204").Value = 45433.38541666666
$ws.Range("G204").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H204").Value = 64.80000305175781
$ws.Range("I204").Value = 45435.38541666666
$ws.Range("I204").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J204").Value = 64.5
$ws.Range("K204").Value = 'High'
$ws.Range("L204").Value = '11/06/2024 05:47:02'

# --- two_line: add rows 38-42 ---
$ws = $wb.Worksheets.Item("two_line")
$ws.Range("A38").Value = 45450.63541666666
$ws.Range("A38").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B38").Value = '11-06-2024 11:15:00'
$ws.Range("C38").Value = 'hour'
$ws.Range("D38").Value = 'SHRIRAMFIN.NS'
$ws.Range("E38").Value = 45446.38541666666
$ws.Range("E38").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F38").Value = 2544.10009765625
$ws.Range("G38").Value = 45449.46875
$ws.Range("G38").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H38").Value = 2535
$ws.Range("I38").Value = 'High'
$ws.Range("J38").Value = '11/06/2024 05:47:02'

$ws.Range("A39").Value = 45453.42708333334
$ws.Range("A39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B39").Value = '11-06-2024 10:15:00'
$ws.Range("C39").Value = 'hour'
$ws.Range("D39").Value = 'ARTSONEN.BO'
$ws.Range("E39").Value = 45439.51041666666
$ws.Range("E39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F39").Value = 189.3999938964844
$ws.Range("G39").Value = 45439.55208333334
$ws.Range("G39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H39").Value = 189.3999938964844
$ws.Range("I39").Value = 'High'
$ws.Range("J39").Value = '11/06/2024 05:47:02'

$ws.Range("A40").Value = 45453.42708333334
$ws.Range("A40").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B40").Value = '11-06-2024 10:15:00'
$ws.Range("C40").Value = 'hour'
$ws.Range("D40").Value = 'ARTSONEN.BO'
$ws.Range("E40").Value = 45439.46875
$ws.Range("E40").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F40").Value = 189.3999938964844
$ws.Range("G40").Value = 45439.55208333334
$ws.Range("G40").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H40").Value = 189.3999938964844
$ws.Range("I40").Value = 'High'
$ws.Range("J40").Value = '11/06/2024 05:47:02'

$ws.Range("A41").Value = 45433.42708333334
$ws.Range("A41").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B41").Value = '11-06-2024 11:15:00'
$ws.Range("C41").Value = 'hour'
$ws.Range("D41").Value = 'GREAVESCOT.NS'
$ws.Range("E41").Value = 45425.38541666666
$ws.Range("E41").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F41").Value = 128.5
$ws.Range("G41").Value = 45425.42708333334
$ws.Range("G41").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H41").Value = 128.5
$ws.Range("I41").Value = 'Low'
$ws.Range("J41").Value = '11/06/2024 05:47:02'

$ws.Range("A42").Value = 45434.55208333334
$ws.Range("A42").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B42").Value = '11-06-2024 10:15:00'
$ws.Range("C42").Value = 'hour'
$ws.Range("D42").Value = 'XPROINDIA.NS'
$ws.Range("E42").Value = 45422.42708333334
$ws.Range("E42").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F42").Value = 1039.699951171875
$ws.Range("G42").Value = 45433.38541666666
$ws.Range("G42").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H42").Value = 1035.650024414062
$ws.Range("I42").Value = 'High'
$ws.Range("J42").Value = '11/06/2024 05:47:02'

# --- ph_pl_breakout_line: add rows 696-698 ---
$ws = $wb.Worksheets.Item("ph_pl_breakout_line")
$ws.Range("A696").Value = 'BGRENERGY.NS'
$ws.Range("B696").Value = 45446.51041666666
$ws.Range("B696").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C696").Value = 40.29999923706055
$ws.Range("D696").Value = 40.29999923706055
$ws.Range("E696").Value = 40.29999923706055
$ws.Range("F696").Value = 'High'
$ws.Range("G696").Value = 40.29999923706055
$ws.Range("H696").Value = 'hour'
$ws.Range("I696").Value = '11-06-2024 10:15:00'
$ws.Range("J696").Value = 40.77000045776367
$ws.Range("K696").Value = 40
$ws.Range("L696").Value = '11/06/2024 05:47:02'

$ws.Range("A697").Value = 'BGRENERGY.NS'
$ws.Range("B697").Value = 45446.55208333334
$ws.Range("B697").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C697").Value = 40.29999923706055
$ws.Range("D697").Value = 40.29999923706055
$ws.Range("E697").Value = 40.29999923706055
$ws.Range("F697").Value = 'High'
$ws.Range("G697").Value = 40.29999923706055
$ws.Range("H697").Value = 'hour'
$ws.Range("I697").Value = '11-06-2024 10:15:00'
$ws.Range("J697").Value = 40.77000045776367
$ws.Range("K697").Value = 40
$ws.Range("L697").Value = '11/06/2024 05:47:02'

$ws.Range("A698").Value = 'BGRENERGY.NS'
$ws.Range("B698").Value = 45446.63541666666
$ws.Range("B698").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C698").Value = 40.29999923706055
$ws.Range("D698").Value = 40.29999923706055
$ws.Range("E698").Value = 40.29999923706055
$ws.Range("F698").Value = 'High'
$ws.Range("G698").Value = 40.29999923706055
$ws.Range("H698").Value = 'hour'
$ws.Range("I698").Value = '11-06-2024 10:15:00'
$ws.Range("J698").Value = 40.77000045776367
$ws.Range("K698").Value = 40
$ws.Range("L698").Value = '11/06/2024 05:47:02'
